$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: "Detailed Test" application
$ws.Cells.Item(7, 1).Value = "Detailed"
$ws.Cells.Item(7, 2).Value = "Test"
$ws.Cells.Item(7, 3).Value = "detailed@example.com"
$ws.Cells.Item(7, 4).Value = "USA"
$ws.Cells.Item(7, 5).Value = "Contributor"
$ws.Cells.Item(7, 6).Value = "Intermediate"
$ws.Cells.Item(7, 7).Value = "Python, JavaScript, React, Node.js"
$ws.Cells.Item(7, 8).Value = "Want to contribute to open source projects and learn from experienced developers"
$ws.Cells.Item(7, 9).Value = "https://github.com/detailedtest"
$ws.Cells.Item(7, 10).Value = "https://detailedtest.dev"
$ws.Cells.Item(7, 11).Value = "10-15 hours per week"
$ws.Cells.Item(7, 12).Value = $true
$ws.Cells.Item(7, 13).Value = $true
$ws.Cells.Item(7, 14).Value = "2025-09-12 22:00:32"

# Row 8: "Complete Test" application
$ws.Cells.Item(8, 1).Value = "Complete"
$ws.Cells.Item(8, 2).Value = "Test"
$ws.Cells.Item(8, 3).Value = "complete@example.com"
$ws.Cells.Item(8, 4).Value = "USA"
$ws.Cells.Item(8, 5).Value = "Contributor"
$ws.Cells.Item(8, 6).Value = "Intermediate"
$ws.Cells.Item(8, 7).Value = "Python, JavaScript, React, Node.js, Excel Integration"
$ws.Cells.Item(8, 8).Value = "Want to test the complete flow from frontend to Excel"
$ws.Cells.Item(8, 9).Value = "https://github.com/completetest"
$ws.Cells.Item(8, 10).Value = "https://completetest.dev"
$ws.Cells.Item(8, 11).Value = "15-20 hours per week"
$ws.Cells.Item(8, 12).Value = $true
$ws.Cells.Item(8, 13).Value = $true
$ws.Cells.Item(8, 14).Value = "2025-09-12 22:08:06"
